$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Cells.Item(1, 4).Value = "desarrolladores"
$ws.Cells.Item(1, 5).Value = "publicadores"
$ws.Cells.Item(1, 6).Value = "fechas_lanzamiento"

# F1 is a brand new header cell; give it the same bold/centered/bordered
# style already used by the other header cells (D1/E1) by copying their format.
$ws.Cells.Item(1, 4).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2..21 -> developers (D), publishers (E), release dates JSON (F)
$rows = @(
    @{ Row = 2;  D = "Lucid Sheep Games";              E = "Lucid Sheep Games";                   F = '{"Japan": "Unreleased", "NorthAmerica": "April 12, 2018", "Europe": "April 12, 2018", "Australia": "April 12, 2018"}' },
    @{ Row = 3;  D = "Beatshapers";                     E = "Beatshapers";                          F = '{"Japan": "Unreleased", "NorthAmerica": "January 24, 2019", "Europe": "January 24, 2019", "Australia": "January 24, 2019"}' },
    @{ Row = 4;  D = "QubicGames";                      E = "QubicGames";                           F = '{"Japan": "Unreleased", "NorthAmerica": "August 9, 2019", "Europe": "August 9, 2019", "Australia": "August 9, 2019"}' },
    @{ Row = 5;  D = "Skobbejak Games";                 E = "Skobbejak Games";                      F = '{"Japan": "Unreleased", "NorthAmerica": "November 6, 2018", "Europe": "December 7, 2018", "Australia": "December 7, 2018"}' },
    @{ Row = 6;  D = "Nintendo EPD";                    E = "Nintendo";                             F = '{"Japan": "March 3, 2017", "NorthAmerica": "March 3, 2017", "Europe": "March 3, 2017", "Australia": "March 3, 2017"}' },
    @{ Row = 7;  D = "Blue Print";                      E = "Blue Print";                           F = '{"Japan": "December 21, 2017", "NorthAmerica": "January 25, 2018", "Europe": "February 22, 2018", "Australia": "February 22, 2018"}' },
    @{ Row = 8;  D = "Nawia Games";                     E = "Nawia Games";                          F = '{"Japan": "Unreleased", "NorthAmerica": "October 25, 2018", "Europe": "October 25, 2018", "Australia": "October 25, 2018"}' },
    @{ Row = 9;  D = "Ink Stains Games";                E = "HypeTrain Digital";                    F = '{"Japan": "Unreleased", "NorthAmerica": "March 5, 2019", "Europe": "March 5, 2019", "Australia": "March 5, 2019"}' },
    @{ Row = 10; D = "Jetdogs, Zoom Out Games";         E = "JetDogs";                              F = '{"Japan": "Unreleased", "NorthAmerica": "March 24, 2019", "Europe": "April 1, 2019", "Australia": "April 1, 2019"}' },
    @{ Row = 11; D = "Roman Uhilg";                     E = "Roman Uhilg";                          F = '{"Japan": "April 26, 2018", "NorthAmerica": "April 27, 2018", "Europe": "April 27, 2018", "Australia": "April 27, 2018"}' },
    @{ Row = 12; D = "Carlsen Games";                   E = "Carlsen Games";                        F = '{"Japan": "Unreleased", "NorthAmerica": "January 9, 2019", "Europe": "January 9, 2019", "Australia": "January 9, 2019"}' },
    @{ Row = 13; D = "Andrade Games";                   E = "Korion";                               F = '{"Japan": "Unreleased", "NorthAmerica": "June 14, 2018", "Europe": "June 14, 2018", "Australia": "June 14, 2018"}' },
    @{ Row = 14; D = "Ink Stories";                     E = "Ink Stories";                          F = '{"Japan": "Unreleased", "NorthAmerica": "August 2, 2018", "Europe": "August 2, 2018", "Australia": "August 2, 2018"}' },
    @{ Row = 15; D = "Baltoro Games";                   E = "Baltoro Games";                        F = '{"Japan": "Unreleased", "NorthAmerica": "September 20, 2019", "Europe": "September 20, 2019", "Australia": "September 20, 2019"}' },
    @{ Row = 16; D = "Millo Games";                     E = "Millo Games";                          F = '{"Japan": "Unreleased", "NorthAmerica": "December 23, 2019", "Europe": "December 23, 2019", "Australia": "December 23, 2019"}' },
    @{ Row = 17; D = "MidBoss";                         E = "MidBoss";                              F = '{"Japan": "December 27, 2018", "NorthAmerica": "August 14, 2018", "Europe": "August 14, 2018", "Australia": "August 14, 2018"}' },
    @{ Row = 18; D = "Batterystaple Games";             E = "Batterystaple Games";                  F = '{"Japan": "July 10, 2018", "NorthAmerica": "July 10, 2018", "Europe": "July 10, 2018", "Australia": "July 10, 2018"}' },
    @{ Row = 19; D = "Digital Bards";                   E = "Digital Bards";                        F = '{"Japan": "February 21, 2019", "NorthAmerica": "May 30, 2019", "Europe": "May 30, 2019", "Australia": "May 30, 2019"}' },
    @{ Row = 20; D = "Petite Games";                    E = "JP: Rainy Frog, WW: Ratalaika Games";  F = '{"Japan": "September 27, 2017", "NorthAmerica": "September 14, 2017", "Europe": "September 14, 2017", "Australia": "September 14, 2017"}' },
    @{ Row = 21; D = "Joindots";                        E = "Joindots";                             F = '{"Japan": "Unreleased", "NorthAmerica": "February 1, 2018", "Europe": "February 1, 2018", "Australia": "February 1, 2018"}' }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
}
